# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.889.75'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '3.309.57'
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'392.89"
$ws.Range("E5").Value = '  -4.55%  '
$ws.Range("D6").Value = "'123.17"
$ws.Range("E6").Value = '  +6.34%  '
$ws.Range("D7").Value = "'0.580"
$ws.Range("E7").Value = '  +1.01%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = "'0.643"
$ws.Range("E9").Value = '  +2.27%  '
$ws.Range("D10").Value = "'0.115"
$ws.Range("E10").Value = '  +0.22%  '
$ws.Range("D11").Value = "'39.84"
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("D13").Value = '3.845.24'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").Value = "'8.06"
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("D15").Value = "'18.86"
$ws.Range("E15").Value = '  -1.57%  '
$ws.Range("D16").Value = '3.331.46'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '60.913.73'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("B18").Value = 'Polygon'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").Value = "'0.983"
$ws.Range("E18").Value = '  -2.34%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = "'10.89"
$ws.Range("E19").Value = '  +0.66%  '
$ws.Range("E20").Value = '  +4.06%  '
$ws.Range("D21").Value = "'3.13"
$ws.Range("E21").Value = '  -6.48%  '
$ws.Range("D22").Value = "'78.31"
$ws.Range("E22").Value = '  +5.44%  '
$ws.Range("D23").Value = "'12.41"
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("D24").Value = "'292.26"
$ws.Range("E24").Value = '  -0.58%  '
$ws.Range("D25").Value = "'3.02"
$ws.Range("E25").Value = '  -3.07%  '
$ws.Range("D26").Value = "'4.66"
$ws.Range("E26").Value = '  +9.90%  '
$ws.Range("D27").Value = "'28.32"
$ws.Range("E27").Value = '  -2.75%  '
$ws.Range("D28").Value = "'7.93"
$ws.Range("E28").Value = '  +4.92%  '
$ws.Range("D29").Value = "'7.29"
$ws.Range("E29").Value = '  -7.07%  '
$ws.Range("D30").Value = "'0.169"
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = "'0.111"
$ws.Range("E32").Value = '  -2.92%  '
$ws.Range("D33").Value = "'11.03"
$ws.Range("E33").Value = '  -3.12%  '
$ws.Range("D34").Value = "'2.46"
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("D35").Value = "'40.24"
$ws.Range("E35").Value = '  -5.60%  '
$ws.Range("D36").Value = "'0.0466"
$ws.Range("E36").Value = '  -4.85%  '
$ws.Range("D37").Value = "'51.72"
$ws.Range("E37").Value = '  -1.35%  '
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("D39").Value = "'3.31"
$ws.Range("E39").Value = '  -4.00%  '
$ws.Range("D40").Value = "'2.83"
$ws.Range("E40").Value = '  -8.59%  '
$ws.Range("D41").Value = "'135.18"
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").Value = "'1.91"
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("D43").Value = "'0.119"
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("D44").Value = "'0.272"
$ws.Range("E44").Value = '  -5.46%  '
$ws.Range("D45").Value = "'16.19"
$ws.Range("E45").Value = '  -1.27%  '
$ws.Range("D46").Value = "'3.75"
$ws.Range("E46").Value = '  -3.02%  '
$ws.Range("E47").Value = '  -1.02%  '
$ws.Range("D48").Value = "'20.71"
$ws.Range("E48").Value = '  -1.99%  '
$ws.Range("D49").Value = '3.642.74'
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("D50").Value = '2.071.90'
$ws.Range("E50").Value = '  -3.88%  '
$ws.Range("D51").Value = "'2.29"
$ws.Range("E51").Value = '  -4.46%  '
